$d = $word.ActiveDocument

$replacements = @(
    @("90÷3=", "97÷2="),
    @("34÷7=", "53÷2="),
    @("66÷4=", "16÷8="),
    @("92÷4=", "71÷7="),
    @("62÷8=", "26÷3="),
    @("20÷2=", "84÷9="),
    @("38÷6=", "31÷6="),
    @("71÷3=", "90÷9="),
    @("51÷6=", "15÷9="),
    @("96÷5=", "14÷7="),
    @("77÷4=", "55÷3="),
    @("40÷4=", "85÷4="),
    @("73÷8=", "28÷4="),
    @("15÷8=", "29÷2="),
    @("89÷8=", "60÷2="),
    @("56÷4=", "88÷7="),
    @("89÷7=", "24÷6="),
    @("78÷6=", "39÷2="),
    @("99÷5=", "27÷3="),
    @("91÷4=", "14÷8="),
    @("51÷3=", "93÷6="),
    @("14÷2=", "39÷3="),
    @("89÷2=", "91÷6="),
    @("21÷7=", "76÷8="),
    @("95÷8=", "20÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
